$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '89.875.94'
$ws.Range("E2").Value = '  -1.24%  '
$ws.Range("D3").Value = '3.081.55'
$ws.Range("E3").Value = '  -2.87%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.24'
$ws.Range("E5").Value = '  +7.95%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '618.21'
$ws.Range("E6").Value = '  -1.41%  '
$ws.Range("E7").Value = '  -11.43%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.356'
$ws.Range("E8").Value = '  -3.87%  '
$ws.Range("E9").Value = '  +0.07%  '
$ws.Range("D10").Value = '3.077.89'
$ws.Range("E10").Value = '  -2.93%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.713'
$ws.Range("E11").Value = '  -6.34%  '
$ws.Range("E12").Value = '  -3.07%  '
$ws.Range("E13").Value = '  -0.60%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.84'
$ws.Range("E14").Value = '  -0.51%  '
$ws.Range("D15").Value = '89.803.27'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.34'
$ws.Range("E16").Value = '  -6.74%  '
$ws.Range("D17").Value = '3.653.73'
$ws.Range("E17").Value = '  -2.82%  '
$ws.Range("D18").Value = '3.084.11'
$ws.Range("E18").Value = '  -3.31%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.78'
$ws.Range("E19").Value = '  +0.30%  '
$ws.Range("E20").Value = '  -1.21%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.77'
$ws.Range("E21").Value = '  -6.27%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '429.59'
$ws.Range("E22").Value = '  -8.88%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.42'
$ws.Range("E23").Value = '  +2.98%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.73'
$ws.Range("E24").Value = '  -4.84%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.56'
$ws.Range("E25").Value = '  -6.34%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '85.97'
$ws.Range("E26").Value = '  -10.46%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.71'
$ws.Range("E27").Value = '  -5.50%  '
$ws.Range("D28").Value = '3.314.08'
$ws.Range("E28").Value = '  -0.88%  '
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.08'
$ws.Range("E30").Value = '  -2.21%  '
$ws.Range("E31").Value = '  +3.17%  '
$ws.Range("E32").Value = '  -4.58%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.191'
$ws.Range("E33").Value = '  -1.42%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '25.46'
$ws.Range("E34").Value = '  -10.38%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.149'
$ws.Range("E35").Value = '  +3.71%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.70'
$ws.Range("E36").Value = '  +1.97%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '495.30'
$ws.Range("E37").Value = '  -5.72%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.00'
$ws.Range("E38").Value = '  +0.21%  '
$ws.Range("E39").Value = '  -3.08%  '
$ws.Range("E40").Value = '  -4.55%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.59'
$ws.Range("E41").Value = '  +53.97%  '
$ws.Range("B42").Value = 'WhiteBITCoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '22.09'
$ws.Range("E42").Value = '  -0.63%  '
$ws.Range("B43").Value = 'Hedera'
$ws.Range("C43").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0860'
$ws.Range("E43").Value = '  -4.67%  '
$ws.Range("E44").Value = '  -0.04%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.397'
$ws.Range("E45").Value = '  -5.59%  '
$ws.Range("E46").Value = '  -6.84%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.673'
$ws.Range("E47").Value = '  -4.66%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '149.31'
$ws.Range("E48").Value = '  -0.61%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '44.38'
$ws.Range("E49").Value = '  -2.33%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.999'
$ws.Range("E50").Value = '  -0.22%  '
$ws.Range("E51").Value = '  -4.96%  '
